$wb = $excel.ActiveWorkbook

$renames = @{
    "summ27911542" = "summ01333973"
    "summ28210844" = "summ01629687"
    "summ28525950" = "summ01916766"
    "summ28904057" = "summ02221265"
    "summ29204330" = "summ02513664"
    "summ29504158" = "summ02812772"
    "summ29820117" = "summ03096818"
    "summ30118976" = "summ03385197"
    "summ30418947" = "summ03671666"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
